$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows before row 401, shifting existing rows 401+ down by 3.
$ws.Rows("401:403").Insert()

# Fill in the new rows (401-403) with data for the new date group (44748).
$newRows = @(
    @{ Row = 401; I = "Especial"; J = 400; K = 11000; L = 12000; M = 11500; P = 639 },
    @{ Row = 402; I = "Primera";  J = 500; K = 9000;  L = 10000; M = 9500;  P = 528 },
    @{ Row = 403; I = "Segunda";  J = 360; K = 6000;  L = 7000;  M = 6500;  P = 361 }
)

foreach ($r in $newRows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value = 2
    $ws.Cells.Item($row, 2).Value = "Comercializadora del Agro de Limarí"
    $ws.Cells.Item($row, 3).Value = "Coquimbo"
    $ws.Cells.Item($row, 4).Value = 44748
    $ws.Cells.Item($row, 5).Value = 4
    $ws.Cells.Item($row, 6).Value = 100112043
    $ws.Cells.Item($row, 7).Value = "Pepino dulce"
    $ws.Cells.Item($row, 8).Value = "Cultivar IV Región"
    $ws.Cells.Item($row, 9).Value = $r.I
    $ws.Cells.Item($row, 10).Value = $r.J
    $ws.Cells.Item($row, 11).Value = $r.K
    $ws.Cells.Item($row, 12).Value = $r.L
    $ws.Cells.Item($row, 13).Value = $r.M
    $ws.Cells.Item($row, 14).Value = "`$/bandeja 18 kilos"
    $ws.Cells.Item($row, 15).Value = "Provincia de Limarí"
    $ws.Cells.Item($row, 16).Value = $r.P
    $ws.Cells.Item($row, 17).Value = 18
    $ws.Cells.Item($row, 18).Value = "Hortaliza"
}

Write-Host "Done"
